$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update B4:B23 ("cliente_id" column) from sequential numbers to cedula-style IDs ---
for ($r = 4; $r -le 23; $r++) {
    $cell = $ws.Cells.Item($r, 2)   # column B
    $old = $cell.Value()
    $cell.Value = 1045012300 + $old
}

# --- 2. Update H4:H53 ("cliente_id" lookup column) the same way ---
for ($r = 4; $r -le 53; $r++) {
    $cell = $ws.Cells.Item($r, 8)   # column H
    $old = $cell.Value()
    $cell.Value = 1045012300 + $old
}

# --- 3. Normalise the border formatting now that both id columns hold the longer
#        cedula-style numbers: every interior row (B5:B22 and H4:H52) adopts the
#        border style already used by B4, and the last row of each table (H53)
#        adopts the style already used by the table's final row (B23).
$ws.Range("B4").Copy()
$ws.Range("B5:B22").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("H4:H52").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B23").Copy()
$ws.Range("H53").PasteSpecial(-4122)      # xlPasteFormats

$excel.CutCopyMode = 0

# --- 4. Widen column H slightly to fit the longer cedula numbers (target stored width: 11) ---
$ws.Columns.Item(8).ColumnWidth = 10.1666666666667

# --- 5. Restore the selection to where the user last clicked ---
$ws.Range("D49").Select()
